$d = $word.ActiveDocument

# The CU17 planilla ("use-case sheet") had its header cell mislabeled as
# "CU18"; correct it back to "CU17" (first cell of the first table).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "CU18",    # FindText
    $true,     # MatchCase
    $true,     # MatchWholeWord
    $false,    # MatchWildcards
    $false,    # MatchSoundsLike
    $false,    # MatchAllWordForms
    $true,     # Forward
    1,         # Wrap (wdFindContinue)
    $false,    # Format
    "CU17",    # ReplaceWith
    2          # Replace (wdReplaceAll)
)
